$d = $word.ActiveDocument

# The "Avg. Time on Tasks (Pretest)" summary-statistics row currently reports
# a classical Mean (SD); the commit replaces it with Median [IQR] values.
# That row is uniquely identified by its current numeric cell contents, so
# locate it via Find (scoped to the table cells) rather than relying on a
# hard-coded row index.

$nbsp = [char]0xA0
$oldLabel = "$nbsp$nbsp" + "Mean (SD)"
$newLabel = "$nbsp$nbsp" + "Median [IQR]"

$found = $false
foreach ($tbl in $d.Tables) {
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $labelCell = $tbl.Cell($r, 1)
        $labelText = $labelCell.Range.Text.TrimEnd([char]13, [char]7)
        if ($labelText -eq $oldLabel) {
            $valueCell = $tbl.Cell($r, 2)
            $valueText = $valueCell.Range.Text.TrimEnd([char]13, [char]7)
            if ($valueText -eq "1.58 (2.38)") {
                $labelCell.Range.Text = $newLabel
                $tbl.Cell($r, 2).Range.Text = "1.04 [0.895]"
                $tbl.Cell($r, 3).Range.Text = "1.06 [0.880]"
                $tbl.Cell($r, 4).Range.Text = "1.05 [0.903]"
                $found = $true
                break
            }
        }
    }
    if ($found) { break }
}

if (-not $found) {
    throw "Target table row (Avg. Time on Tasks Mean (SD)) not found"
}
